# Weekly Time Record - fill in Monday's "Out" time and Tuesday's In/Out times,
# move the active selection down to D15, re-stamp the (duplicated) print-area
# defined name twice more (artifact of re-setting the print area), and nudge
# the column widths that Calc recomputed when the page layout was touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Time Record")

# --- Time entries -----------------------------------------------------
# Monday (row 13): clocked out at 8:30 PM (was blank)
$ws.Range("D13").Value = 0.854166666666667

# Tuesday (row 14): clocked in at 11:30 AM, out at 3:30 PM (both were blank)
$ws.Range("C14").Value = 0.479166666666667
$ws.Range("D14").Value = 0.645833333333333

# --- Selected cell moved from D13 to D15 -------------------------------
$ws.Range("D15").Select() | Out-Null

# --- Re-assert the print area (Calc appends a fresh _xlnm.Print_Area_* --
# --- defined name each time this is done instead of reusing one) -------
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

# --- Column widths recomputed by the same print-area refresh -----------
$ws.Columns.Item(2).ColumnWidth = 8.83
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 6.83
$ws.Columns.Item(8).ColumnWidth = 8.0
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 6.83
